$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure column D cells keep their original text/string representation
# (prices like "1.001" or "10.00" would otherwise be auto-converted to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.780.19'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').Value = '1.612.48'
$ws.Range('E3').Value = '  -3.72%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '208.59'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = '0.5180'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '0.2560'
$ws.Range('E8').Value = '  -3.51%  '
$ws.Range('D9').Value = '0.06200'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').Value = '20.25'
$ws.Range('E10').Value = '  -5.02%  '
$ws.Range('D11').Value = '0.07523'
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = '1.634.51'
$ws.Range('E12').Value = '  -2.51%  '
$ws.Range('D13').Value = '4.351'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').Value = '1.832.53'
$ws.Range('E14').Value = '  -3.79%  '
$ws.Range('D15').Value = '0.5411'
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').Value = '0.0₅7869'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = '63.76'
$ws.Range('E17').Value = '  -4.80%  '
$ws.Range('D18').Value = '25.772.46'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '4.612'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('D21').Value = '182.97'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').Value = '10.00'
$ws.Range('E22').Value = '  -4.06%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '6.035'
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '144.38'
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').Value = '0.1197'
$ws.Range('E26').Value = '  -4.60%  '
$ws.Range('D27').Value = '7.333'
$ws.Range('E27').Value = '  -3.31%  '
$ws.Range('D28').Value = '15.43'
$ws.Range('E28').Value = '  -3.26%  '
$ws.Range('D29').Value = '1.351'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('D30').Value = '0.05851'
$ws.Range('E30').Value = '  -6.07%  '
$ws.Range('D31').Value = '1.235'
$ws.Range('E31').Value = '  -3.90%  '
$ws.Range('D32').Value = '3.360'
$ws.Range('E32').Value = '  -4.29%  '
$ws.Range('D33').Value = '3.326'
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('D34').Value = '1.594'
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('D35').Value = '0.9641'
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('D36').Value = '2.381'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('D37').Value = '2.706'
$ws.Range('D38').Value = '0.5727'
$ws.Range('E38').Value = '  -5.60%  '
$ws.Range('D39').Value = '0.01579'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').Value = '0.8387'
$ws.Range('E41').Value = '  -3.79%  '
$ws.Range('D42').Value = '5.644'
$ws.Range('E42').Value = '  -7.72%  '
$ws.Range('D43').Value = '1.016.49'
$ws.Range('E43').Value = '  -7.72%  '
$ws.Range('D44').Value = '98.96'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').Value = '1.759.28'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('D47').Value = '1.007'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '54.07'
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('D49').Value = '7.913'
$ws.Range('E49').Value = '  -1.65%  '
$ws.Range('D50').Value = '0.05150'
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('D51').Value = '0.4204'
$ws.Range('E51').Value = '  -1.22%  '
